$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "total" column sub-header cell (was "unnamed: 1_level_1")
$ws.Range("B2").Value = "total"

# Remove the two section-header-only rows that carried no numeric data:
# row 5 ("situacao do domicilio") and, after that row is gone, the new row 7
# ("grandes regioes e unidades da federacao"). Everything below shifts up.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(7).Delete()

# Corrected/updated numeric data (B:I) for every remaining data row (4-38)
$data = @(
  @(4, 0.45, 0.8100000000000001, 0.54, 0.65, 0.67, 0.53, 0.92, 0.85),
  @(5, 0.47, 0.91, 0.58, 0.7, 0.72, 0.5600000000000001, 0.95, 0.87),
  @(6, 1.08, 1.28, 1.08, 1.46, 1.88, 1.71, 3.58, 2.41),
  @(7, 2.64, 4.07, 3.12, 2.96, 3.35, 2.23, 2.89, 2.54),
  @(8, 2.06, 2.83, 2.58, 3.28, 3.7, 2.76, 5.67, 4.51),
  @(9, 5.24, 6.58, 6.11, 5.9, 8.710000000000001, 7.4, 6.15, 6.48),
  @(10, 3.11, 4.79, 3.43, 4.09, 5.15, 3.5, 5.47, 4.68),
  @(11, 3.07, 5.05, 3.2, 8.279999999999999, 7.09, 3.87, 7.87, 6.24),
  @(12, 6.28, 8.630000000000001, 6.72, 7.01, 8.07, 5.62, 7.74, 7.21),
  @(13, 7.45, 13.21, 9.779999999999999, 8.859999999999999, 9.58, 5.55, 12.5, 8.119999999999999),
  @(14, 3.22, 6.28, 3.38, 4.22, 3.67, 2.43, 4.19, 3.24),
  @(15, 0.98, 1.31, 1.16, 1.43, 1.31, 1.19, 2.36, 2.13),
  @(16, 3.53, 3.81, 4.39, 7.48, 3.44, 3.99, 10.2, 8.43),
  @(17, 2.65, 2.1, 2.43, 3.59, 5.3, 4.12, 10.91, 7.35),
  @(18, 1.76, 3.03, 2.1, 2.36, 2.59, 2.26, 6.07, 4.29),
  @(19, 6.6, 7.37, 7.53, 7.52, 7.34, 7.94, 9.220000000000001, 9.359999999999999),
  @(20, 5.2, 7.45, 6.12, 4.35, 5.19, 4.59, 8.67, 9.390000000000001),
  @(21, 1.77, 2.46, 2.02, 2.32, 3.69, 2.45, 4.39, 4.57),
  @(22, 4.7, 4.92, 5.7, 6.46, 7.98, 5.56, 10.73, 8.68),
  @(23, 4.32, 4.39, 5.32, 5.13, 7.95, 4.92, 7.48, 6.48),
  @(24, 1.83, 2.75, 1.99, 2.28, 2.13, 2.37, 4.54, 3.67),
  @(25, 0.63, 1.22, 0.73, 1.04, 1.08, 0.8, 1.44, 1.29),
  @(26, 1.3, 2.06, 1.4, 1.86, 1.92, 1.68, 2.65, 2.24),
  @(27, 2.41, 3.55, 2.52, 3.52, 5.32, 4.13, 7.2, 5.52),
  @(28, 1.73, 3.27, 1.87, 2.37, 2.58, 2.08, 3.5, 2.99),
  @(29, 0.8100000000000001, 1.8, 1, 1.54, 1.55, 1.04, 1.94, 1.77),
  @(30, 0.78, 1.54, 0.92, 1.17, 1.27, 1.04, 1.72, 1.49),
  @(31, 1.31, 2.58, 1.75, 2.18, 2.15, 1.57, 3.01, 2.88),
  @(32, 1.79, 2.44, 1.95, 2.44, 3.21, 2.33, 3.37, 2.43),
  @(33, 1.1, 2.29, 1.25, 1.54, 1.76, 1.58, 2.43, 2.08),
  @(34, 1.13, 2.1, 1.35, 1.87, 1.98, 1.47, 2.23, 1.85),
  @(35, 2.08, 4.28, 2.74, 2.93, 4.08, 2.56, 4.01, 4.05),
  @(36, 2.44, 4.98, 2.82, 3.91, 4.34, 3.11, 6.4, 6.63),
  @(37, 1.79, 2.82, 1.94, 3.15, 2.81, 2.3, 3.48, 3.21),
  @(38, 3.12, 6.21, 4.18, 4.94, 5.5, 3.92, 4.37, 2.49)
)

foreach ($rowspec in $data) {
  $r = $rowspec[0]
  for ($c = 0; $c -lt 8; $c++) {
    $ws.Cells.Item($r, $c + 2).Value = $rowspec[$c + 1]
  }
}

Write-Output "edit applied"
